# Add Switzerland ("Swiss") market test-data sheet, based on the existing
# "Czech" sheet layout, and make it the active sheet.

$wb = $excel.ActiveWorkbook
$czech = $wb.Worksheets.Item("Czech")

# Duplicate the Czech sheet (keeps formatting, column widths, merged cells,
# styles, etc.) and place the copy immediately after it.
$czech.Copy($null, $czech)
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# Update the market name and Jira reference for Switzerland.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2352"

# Insert the extra "PR1DSCH" repeater row after the existing "PR1DS" row.
$swiss.Rows.Item(9).Insert()
$swiss.Range("A8").Copy()
$swiss.Range("A9").PasteSpecial(-4122)   # xlPasteFormats
$swiss.Range("A9").Value = "PR1DSCH"

# Insert the extra "PR8ASCH" repeater row after the existing "PR8AS" row.
$swiss.Rows.Item(11).Insert()
$swiss.Range("A10").Copy()
$swiss.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$swiss.Range("A11").Value = "PR8ASCH"

# The Czech sheet is no longer the active tab, so clear its special
# "current selection" highlight back to a plain full-sheet selection.
$czech.Cells.Select()

# Make the new Swiss sheet the active tab, matching where the workbook was
# last left selected.
$swiss.Range("A8").Select()
$swiss.Activate()
